$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D15").Formula = "=25.9/4+12/4"
$ws.Range("F15").Formula = "=64800/4+45200/4"

$ws.Range("E19").Select()
